$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11: PORTB / 50~53 and 10~13
$ws.Range("A11").Value = "PORTB"
$ws.Range("B11").Value = "50~53 and 10~13"

# Change B6 value from "10~12" to "6~8"
$ws.Range("B6").Value = "6~8"

# Set column B width to approximate the bestFit width seen in diff (15.28515625)
$ws.Columns.Item(2).ColumnWidth = 14.5

# Update active selection to B7 as in the diff
$ws.Range("B7").Select() | Out-Null
